$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rename specific BoM items to their "_V2" revision.
$ws.Range("C32").Value = "Top_Scale_V2"
$ws.Range("C33").Value = "Top_Ring_Big_Inside_V2"
$ws.Range("C34").Value = "Top_Ring_Smal_Inside_V2"
$ws.Range("C35").Value = "Top_Ring_Big_Outer_V2"
$ws.Range("C36").Value = "Top_Ring_Small_Inside_V2"
$ws.Range("C40").Value = "Wall_Inside_V2"
$ws.Range("C48").Value = "Wall_Outside_V2"

# Move the active selection the same way it ended up in the saved file.
$ws.Range("C48").Select()
